# Update gh-pages to output generated at 456a3b4
$wb = $excel.ActiveWorkbook

# Sheet "展览" (exhibitions)
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 826
$ws1.Range("F3").Value = 27
$ws1.Range("F9").Value = 552
$ws1.Range("F13").Value = 13486
$ws1.Range("F17").Value = 5554
$ws1.Range("F18").Value = 5581
$ws1.Range("F19").Value = 56

# Sheet "演出" (performances)
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F24").Value = 10

# Sheet "全部类型" (all types)
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 826
$ws4.Range("F10").Value = 27
$ws4.Range("F31").Value = 552
$ws4.Range("F35").Value = 13486
$ws4.Range("F38").Value = 10
$ws4.Range("F40").Value = 5554
$ws4.Range("F41").Value = 5581
$ws4.Range("F42").Value = 56
